# Publicación y S Postgrado
# 1- Publicación SMQ
# 2- Tutorización Postgrado 2023-2024
#
# Replace the two existing supervision rows (row 2: MSc in Psychology /
# Yenny Johanna Baron Londoño; row 3: MSc in Neuropsychology / Sara Silva
# Gómez) with three rows of "MSc in Neuropsychology" / "2023-2024" /
# Universidad Internacional de Valencia supervisions, one per student.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$what  = "MSc in Neuropsychology"
$when  = "2023-2024"
$where = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"

$students = @(
    "Leidy Nathaly Peláez Bernal",
    "Jimena Zanizo Chambi",
    "Liceth Andrea Zaraza Osorio"
)

for ($i = 0; $i -lt $students.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $what
    $ws.Cells.Item($row, 2).Value = $when
    $ws.Cells.Item($row, 3).Value = $students[$i]
    $ws.Cells.Item($row, 4).Value = $where
}
$ws.Rows.Item(4).RowHeight = 57.6

$ws.Range("C6").Select()
